$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the Cash Flow Statement section whose D:J values shift one column to
# the left (D<-old E, E<-old F, ... I<-old J), with a brand new figure (or "NA")
# placed into the now-vacated column J.
# New column J values keyed by row number; "NA" rows get the text marker,
# others get a freshly reported numeric figure.
$newJValues = @{
    83  = "NA"
    89  = -10200
    91  = -300
    94  = "NA"
    100 = "NA"
    101 = "NA"
    102 = -5700
}

foreach ($r in 83, 89, 91, 94, 100, 101, 102) {
    # Capture existing E:J values before overwriting anything.
    $vals = @()
    for ($c = 5; $c -le 10; $c++) {
        $vals += $ws.Cells.Item($r, $c).Value2
    }

    # Shift them left into D:I.
    for ($c = 4; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 4]
    }

    # Write the new value (number or "NA" text) into column J.
    $ws.Cells.Item($r, 10).Value = $newJValues[$r]
}
